# Iteración_2 Registro.xlsx - "Descripción de CU - 24  Diagramas de secuencia y robustez"
#
# Updates progress-tracking numbers on the "Casos de Uso" sheet for a handful
# of tasks (rows 26, 31, 33-35, 42-43) and records a brand-new completed task
# row (44: "CU - 24 / Editar Perfil / Descripción") that mirrors the style of
# the rows directly above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Row 26 (CU - 21 / Administrar cliente): record 0.5h consumed in the 9th day slot
$ws.Range("AI26").Value = 0.5

# --- Row 31 (CU - 06 / Registrar pago colaborador): move the day-11 consumption
#     registered previously (0.25h, col AR) to day-10 (col AC = 0.5h) and add a
#     fresh 1h entry on day-12 (col AL)
$ws.Range("AC31").Value = 0.5
$ws.Range("AL31").Value = 1
$ws.Range("AR31").Value = ""

# --- Rows 33, 34, 35, 42, 43: status column F flips from "Por iniciar" to "Hecho"
$ws.Range("F33").Value = "Hecho"
$ws.Range("F34").Value = "Hecho"
$ws.Range("F35").Value = "Hecho"
$ws.Range("F42").Value = "Hecho"
$ws.Range("F43").Value = "Hecho"

# --- Row 34 (Cobrar mensualidad): 0.5h consumed on day-13 (col AO)
$ws.Range("AO34").Value = 0.5

# --- Row 35 (Administrar cliente): 1h consumed on day-13 (col AO)
$ws.Range("AO35").Value = 1

# --- Row 42 (CU -24 / Editar Perfil / Diagrama de Robustez): 0.5h consumed on day-9 (col AI)
$ws.Range("AI42").Value = 0.5

# --- Row 43 (CU - 24 / Editar Perfil / Diagrama de Secuencia): 0.7h consumed on day-9 (col AI)
$ws.Range("AI43").Value = 0.7

# --- Row 44: brand-new task line (CU - 24 / Editar Perfil / Descripción / Mauricio / Hecho)
#     First copy the visual formatting from the row above (43) so the new row gets
#     the same "done" shading/style as its neighbours, then fill in the content.
$ws.Range("B43:BA43").Copy()
$ws.Range("B44:BA44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B44").Value = "CU - 24"
$ws.Range("C44").Value = "Editar Perfil"
$ws.Range("D44").Value = "Descripción"
$ws.Range("E44").Value = "Mauricio"
$ws.Range("F44").Value = "Hecho"
$ws.Range("G44").Value = 0.5
$ws.Range("AC44").Value = 0.25
$ws.Range("AU44").Value = 0.15

# --- Selection / view bookkeeping (cosmetic - matches the author's last-saved cursor position)
$ws.Range("F37").Select()
